$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @(
    "FT2316800GH03LFZ",
    "FT231680T6161WYX",
    "FT231680T6162257",
    "FT2316800GH04LW6",
    "FT2316805D7D108X"
)

$startRow = 11
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}
